$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 500
$ws.Range("I31").Value = 500
$ws.Range("K31").Value = 1500
$ws.Range("M31").Value = -1270

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 545.1667
$ws.Range("I33").Value = 700.5
$ws.Range("J33").Value = 467.5
$ws.Range("K33").Value = 700.5
$ws.Range("L33").Value = 467.5
$ws.Range("M33").Value = -471.5
$ws.Range("N33").Value = -925.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4160.1816
$ws.Range("I64").Value = 3800
$ws.Range("J64").Value = 4409.5386
$ws.Range("K64").Value = 3800
$ws.Range("L64").Value = 4409.5386
$ws.Range("M64").Value = -3552
$ws.Range("N64").Value = -4905.5386

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 4160.1816
$ws.Range("I67").Value = 3800
$ws.Range("J67").Value = 4409.5386
$ws.Range("K67").Value = 3800
$ws.Range("L67").Value = 4409.5386
$ws.Range("M67").Value = -2942
$ws.Range("N67").Value = -6125.5386

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3774.2173
$ws.Range("I74").Value = 3653.5334
$ws.Range("J74").Value = 4000.5
$ws.Range("K74").Value = 3653.5334
$ws.Range("L74").Value = 4000.5
$ws.Range("M74").Value = -2717.5334
$ws.Range("N74").Value = -5872.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3774.2173
$ws.Range("I77").Value = 3653.5334
$ws.Range("J77").Value = 4000.5
$ws.Range("K77").Value = 18267.667
$ws.Range("L77").Value = 20002.5
$ws.Range("M77").Value = -13587.667
$ws.Range("N77").Value = -29362.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 39800
$ws.Range("J87").Value = 39800
$ws.Range("L87").Value = 39800
$ws.Range("N87").Value = -42296

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 39800
$ws.Range("J90").Value = 39800
$ws.Range("L90").Value = 119400
$ws.Range("N90").Value = -131880

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1890255.6
$ws.Range("I137").Value = 2778798.8
$ws.Range("J137").Value = 8634.823
$ws.Range("K137").Value = 8336396.399999999
$ws.Range("L137").Value = 25904.469
$ws.Range("M137").Value = -8333846.399999999
$ws.Range("N137").Value = -31004.469

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1487.826
$ws.Range("I2").Value = 1487.7273
$ws.Range("J2").Value = 1488.0769
$ws.Range("K2").Value = 1487.7273
$ws.Range("L2").Value = 1488.0769
$ws.Range("M2").Value = -1374.7273
$ws.Range("N2").Value = -1714.0769

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1609.75
$ws.Range("I45").Value = 1594.75
$ws.Range("J45").Value = 1624.75
$ws.Range("K45").Value = 1594.75
$ws.Range("L45").Value = 1624.75
$ws.Range("M45").Value = -1217.75
$ws.Range("N45").Value = -2378.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 17894978
$ws.Range("I61").Value = 20430318
$ws.Range("J61").Value = 147602
$ws.Range("K61").Value = 20430318
$ws.Range("L61").Value = 147602
$ws.Range("M61").Value = -20430106
$ws.Range("N61").Value = -148026

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3225.5
$ws.Range("I63").Value = 2811
$ws.Range("J63").Value = 3640
$ws.Range("K63").Value = 2811
$ws.Range("L63").Value = 3640
$ws.Range("M63").Value = -2125
$ws.Range("N63").Value = -5012

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3225.5
$ws.Range("I66").Value = 2811
$ws.Range("J66").Value = 3640
$ws.Range("K66").Value = 14055
$ws.Range("L66").Value = 18200
$ws.Range("M66").Value = -10623
$ws.Range("N66").Value = -25064

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 10480737
$ws.Range("I74").Value = 15198108
$ws.Range("J74").Value = 102520.2
$ws.Range("K74").Value = 15198108
$ws.Range("L74").Value = 102520.2
$ws.Range("M74").Value = -15197234
$ws.Range("N74").Value = -104268.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 10480737
$ws.Range("I77").Value = 15198108
$ws.Range("J77").Value = 102520.2
$ws.Range("K77").Value = 75990540
$ws.Range("L77").Value = 512601
$ws.Range("M77").Value = -75986172
$ws.Range("N77").Value = -521337

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 4167628.5
$ws.Range("I97").Value = 5209387.5
$ws.Range("K97").Value = 5209387.5
$ws.Range("M97").Value = -5208891.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1487.826
$ws.Range("I116").Value = 1487.7273
$ws.Range("J116").Value = 1488.0769
$ws.Range("K116").Value = 1487.7273
$ws.Range("L116").Value = 1488.0769
$ws.Range("M116").Value = 806.2727
$ws.Range("N116").Value = -6076.0769

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 45464.105
$ws.Range("I132").Value = 30239.4
$ws.Range("J132").Value = 86453.69500000001
$ws.Range("K132").Value = 90718.20000000001
$ws.Range("L132").Value = 259361.085
$ws.Range("M132").Value = -88188.20000000001
$ws.Range("N132").Value = -264421.085

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 17894978
$ws.Range("I136").Value = 20430318
$ws.Range("J136").Value = 147602
$ws.Range("K136").Value = 61290954
$ws.Range("L136").Value = 442806
$ws.Range("M136").Value = -61288404
$ws.Range("N136").Value = -447906

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1487.826
$ws.Range("I3").Value = 1487.7273
$ws.Range("J3").Value = 1488.0769
$ws.Range("K3").Value = 1487.7273
$ws.Range("L3").Value = 1488.0769
$ws.Range("M3").Value = -1373.7273
$ws.Range("N3").Value = -1716.0769

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1488
$ws.Range("I20").Value = 1336
$ws.Range("J20").Value = 1533.6
$ws.Range("K20").Value = 1336
$ws.Range("L20").Value = 1533.6
$ws.Range("M20").Value = -1089
$ws.Range("N20").Value = -2027.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1521.1666
$ws.Range("I94").Value = 1125.4
$ws.Range("J94").Value = 3500
$ws.Range("K94").Value = 1125.4
$ws.Range("L94").Value = 3500
$ws.Range("M94").Value = -674.4000000000001
$ws.Range("N94").Value = -4402

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1120.0667
$ws.Range("I99").Value = 999
$ws.Range("K99").Value = 999
$ws.Range("M99").Value = 499

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 41668144
$ws.Range("I105").Value = 45456020
$ws.Range("K105").Value = 45456020
$ws.Range("M105").Value = -45454273

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1860.4468
$ws.Range("I134").Value = 1836.4878
$ws.Range("J134").Value = 2024.1666
$ws.Range("K134").Value = 5509.463400000001
$ws.Range("L134").Value = 6072.4998
$ws.Range("M134").Value = -2974.463400000001
$ws.Range("N134").Value = -11142.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 2605.8518
$ws.Range("I103").Value = 413.3
$ws.Range("J103").Value = 3895.5881
$ws.Range("K103").Value = 1239.9
$ws.Range("L103").Value = 11686.7643
$ws.Range("M103").Value = -360.9000000000001
$ws.Range("N103").Value = -13444.7643

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value = 3061.1875
$ws.Range("J115").Value = 3065.2666
$ws.Range("L115").Value = 9195.799800000001
$ws.Range("N115").Value = -11545.7998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 853.4783
$ws.Range("J131").Value = 894.14636
$ws.Range("L131").Value = 2682.43908
$ws.Range("N131").Value = -12762.43908

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1962.5834
$ws.Range("I113").Value = 1265.7273
$ws.Range("J113").Value = 2552.2307
$ws.Range("K113").Value = 1265.7273
$ws.Range("L113").Value = 2552.2307
$ws.Range("M113").Value = 904.2727
$ws.Range("N113").Value = -6892.2307

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 50027.805
$ws.Range("I132").Value = 29619.371
$ws.Range("K132").Value = 88858.113
$ws.Range("M132").Value = -86328.113

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 949.8
$ws.Range("I93").Value = 977.55554
$ws.Range("J93").Value = 700
$ws.Range("K93").Value = 977.55554
$ws.Range("L93").Value = 700
$ws.Range("M93").Value = 270.44446
$ws.Range("N93").Value = -3196

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 30633.334
$ws.Range("J94").Value = 30633.334
$ws.Range("L94").Value = 30633.334
$ws.Range("N94").Value = -31985.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2280.3125
$ws.Range("I81").Value = 795
$ws.Range("J81").Value = 2775.4167
$ws.Range("K81").Value = 1590
$ws.Range("L81").Value = 5550.8334
$ws.Range("M81").Value = -529
$ws.Range("N81").Value = -7672.8334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 2280.3125
$ws.Range("I84").Value = 795
$ws.Range("J84").Value = 2775.4167
$ws.Range("K84").Value = 7950
$ws.Range("L84").Value = 27754.167
$ws.Range("M84").Value = -2646
$ws.Range("N84").Value = -38362.167

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 59876.707
$ws.Range("I100").Value = 50510.5
$ws.Range("J100").Value = 73257
$ws.Range("K100").Value = 101021
$ws.Range("L100").Value = 146514
$ws.Range("M100").Value = -100480
$ws.Range("N100").Value = -147596

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H130").Value = 41309.668
$ws.Range("J130").Value = 41309.668
$ws.Range("L130").Value = 41309.668
$ws.Range("N130").Value = -51349.668
